# Apply updated "想去人数" (F column) counts across the workbook's sheets.
# Source data refresh: gh-pages output generated at 456a3b4

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 316
$ws.Range("F6").Value = 323
$ws.Range("F7").Value = 1133
$ws.Range("F8").Value = 435
$ws.Range("F9").Value = 6989
$ws.Range("F13").Value = 7886
$ws.Range("F16").Value = 5466
$ws.Range("F18").Value = 2346
$ws.Range("F19").Value = 998
$ws.Range("F21").Value = 280
$ws.Range("F24").Value = 7
$ws.Range("F25").Value = 336
$ws.Range("F28").Value = 2164
$ws.Range("F30").Value = 246
$ws.Range("F32").Value = 79
$ws.Range("F33").Value = 555
$ws.Range("F38").Value = 1
$ws.Range("F39").Value = 2195
$ws.Range("F40").Value = 2193

# --- Sheet "演出" (performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 47

# --- Sheet "本地生活" (local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1269

# --- Sheet "全部类型" (all types, aggregated) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1269
$ws.Range("F5").Value = 316
$ws.Range("F9").Value = 323
$ws.Range("F10").Value = 1133
$ws.Range("F11").Value = 435
$ws.Range("F12").Value = 6989
$ws.Range("F16").Value = 7886
$ws.Range("F19").Value = 5466
$ws.Range("F21").Value = 2346
$ws.Range("F22").Value = 998
$ws.Range("F24").Value = 280
$ws.Range("F28").Value = 7
$ws.Range("F29").Value = 47
$ws.Range("F30").Value = 336
$ws.Range("F33").Value = 2164
$ws.Range("F35").Value = 246
$ws.Range("F37").Value = 79
$ws.Range("F38").Value = 555
$ws.Range("F44").Value = 1
$ws.Range("F45").Value = 2195
$ws.Range("F47").Value = 2193
